$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for all existing
#    data rows (2-409) from 45204 to 45205.
$ws.Range("C2:C409").Value = 45205

# 2. Insert three new blank rows before the current last data row (410),
#    pushing the existing row 410 (A 47593-2023 / Sveaskog) down to row 413.
$ws.Range("A410:A412").EntireRow.Insert()

# 3. Fill in the new row 410: A 47664-2023
$ws.Range("A410").Value = "A 47664-2023"
$ws.Range("B410").Value = 45203
$ws.Range("B410").NumberFormat = "YYYY-MM-DD"
$ws.Range("C410").Value = 45205
$ws.Range("C410").NumberFormat = "YYYY-MM-DD"
$ws.Range("D410").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E410").Value = "MALÅ"
$ws.Range("G410").Value = 5.4
$ws.Range("H410:Q410").Value = 0
$ws.Range("R410").WrapText = $true
$ws.Rows.Item(410).RowHeight = 15

# 4. Fill in the new row 411: A 47665-2023
$ws.Range("A411").Value = "A 47665-2023"
$ws.Range("B411").Value = 45203
$ws.Range("B411").NumberFormat = "YYYY-MM-DD"
$ws.Range("C411").Value = 45205
$ws.Range("C411").NumberFormat = "YYYY-MM-DD"
$ws.Range("D411").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E411").Value = "MALÅ"
$ws.Range("G411").Value = 11.8
$ws.Range("H411:Q411").Value = 0
$ws.Range("R411").WrapText = $true
$ws.Rows.Item(411).RowHeight = 15

# 5. Fill in the new row 412: A 47687-2023
$ws.Range("A412").Value = "A 47687-2023"
$ws.Range("B412").Value = 45203
$ws.Range("B412").NumberFormat = "YYYY-MM-DD"
$ws.Range("C412").Value = 45205
$ws.Range("C412").NumberFormat = "YYYY-MM-DD"
$ws.Range("D412").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E412").Value = "MALÅ"
$ws.Range("G412").Value = 56.9
$ws.Range("H412:Q412").Value = 0
$ws.Range("R412").WrapText = $true
$ws.Rows.Item(412).RowHeight = 15

# 6. Row 413 already holds the original row-410 data (A 47593-2023,
#    Sveaskog, 23.1 ha) after the insert shifted it down. Only the
#    "Förändrad" date needs to change, and it picks up the same row
#    height treatment as its new siblings.
$ws.Range("C413").Value = 45205
$ws.Rows.Item(413).RowHeight = 15

# 7. Append the brand new final row 414: A 47666-2023
$ws.Range("A414").Value = "A 47666-2023"
$ws.Range("B414").Value = 45203
$ws.Range("B414").NumberFormat = "YYYY-MM-DD"
$ws.Range("C414").Value = 45205
$ws.Range("C414").NumberFormat = "YYYY-MM-DD"
$ws.Range("D414").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E414").Value = "MALÅ"
$ws.Range("G414").Value = 6.7
$ws.Range("H414:Q414").Value = 0
$ws.Range("R414").WrapText = $true
